$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Bongo-Items")
$ws.Range('B2').Value = '779D81F6-D578-4E22-AAAE-D6118CD9920F'
$ws.Range('E2').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('F2').Value = 'Gold'
$ws.Range('G2').Value = 'Lemon'
$ws.Range('B3').Value = 'D6FF78FE-753E-46A3-8B86-579FA87BB360'
$ws.Range('E3').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('F3').Value = 'Blue'
$ws.Range('G3').Value = 'Coffee'
$ws.Range('B4').Value = '4E0191E4-A89D-4108-A1B0-5FF0BF2DE6D2'
$ws.Range('E4').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('F4').Value = 'Brown'
$ws.Range('G4').Value = 'Raspberry'
$ws.Range('B5').Value = '5E81D482-0EAE-4373-8B0E-48F5533D62C3'
$ws.Range('E5').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('F5').Value = 'Red'
$ws.Range('G5').Value = 'Grape'

$ws = $wb.Worksheets.Item("Bongo-Item-Images")
$ws.Range('B2').Value = '779D81F6-D578-4E22-AAAE-D6118CD9920F'
$ws.Range('C2').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D2').Value = '../images/raccoon.jpeg'
$ws.Range('B3').Value = 'D6FF78FE-753E-46A3-8B86-579FA87BB360'
$ws.Range('C3').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D3').Value = '../images/runningdog.jpeg'
$ws.Range('B4').Value = '4E0191E4-A89D-4108-A1B0-5FF0BF2DE6D2'
$ws.Range('C4').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D4').Value = '../iamges/apple.jpeg'
$ws.Range('B5').Value = '5E81D482-0EAE-4373-8B0E-48F5533D62C3'
$ws.Range('C5').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D5').Value = '../images/broccoli.jpeg'

$ws = $wb.Worksheets.Item("Bongo-Tests")
$ws.Range('B2').Value = '779D81F6-D578-4E22-AAAE-D6118CD9920F'
$ws.Range('C2').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D2').Value = 0
$ws.Range('B3').Value = 'D6FF78FE-753E-46A3-8B86-579FA87BB360'
$ws.Range('C3').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D3').Value = 80
$ws.Range('E3').Value = 'Pass'
$ws.Range('B4').Value = '4E0191E4-A89D-4108-A1B0-5FF0BF2DE6D2'
$ws.Range('C4').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D4').Value = 33
$ws.Range('B5').Value = '5E81D482-0EAE-4373-8B0E-48F5533D62C3'
$ws.Range('C5').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D5').Value = 68

$ws = $wb.Worksheets.Item("Bongo-Test-Images")
$ws.Range('B2').Value = '779D81F6-D578-4E22-AAAE-D6118CD9920F'
$ws.Range('C2').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D2').Value = '../images/dice.jpeg'
$ws.Range('B3').Value = 'D6FF78FE-753E-46A3-8B86-579FA87BB360'
$ws.Range('C3').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D3').Value = '../images/runningdog.jpeg'
$ws.Range('B4').Value = '4E0191E4-A89D-4108-A1B0-5FF0BF2DE6D2'
$ws.Range('C4').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D4').Value = '../images/dice.jpeg'
$ws.Range('B5').Value = '5E81D482-0EAE-4373-8B0E-48F5533D62C3'
$ws.Range('C5').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D5').Value = '../iamges/apple.jpeg'

$ws = $wb.Worksheets.Item("Biff-Items")
$ws.Range('B2').Value = '4644A72F-3224-440D-8657-1825AAC2440C'
$ws.Range('E2').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('F2').Value = 86.28
$ws.Range('G2').Value = 105.02
$ws.Range('H2').Value = 115.67
$ws.Range('I2').Value = '779D81F6-D578-4E22-AAAE-D6118CD9920F'
$ws.Range('J2').Value = 'D6FF78FE-753E-46A3-8B86-579FA87BB360'
$ws.Range('B3').Value = '6A713058-863E-45C4-B05B-CF5671E33F0C'
$ws.Range('E3').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('F3').Value = 110.77
$ws.Range('G3').Value = 83.31999999999999
$ws.Range('H3').Value = 116.18
$ws.Range('I3').Value = '4E0191E4-A89D-4108-A1B0-5FF0BF2DE6D2'
$ws.Range('J3').Value = '5E81D482-0EAE-4373-8B0E-48F5533D62C3'

$ws = $wb.Worksheets.Item("Biff-Tests")
$ws.Range('B2').Value = '4644A72F-3224-440D-8657-1825AAC2440C'
$ws.Range('C2').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D2').Value = 64
$ws.Range('B3').Value = '6A713058-863E-45C4-B05B-CF5671E33F0C'
$ws.Range('C3').Value = 'generated 2023-10-12 09:12:50'
$ws.Range('D3').Value = 77
